# Update the BC-LEEP archetype's basement wall height ("Opt_Bsm_Height") in
# the "Geometry options" lookup table. This single input drives (via VLOOKUP
# and a long chain of dependent formulas on "geometry calculation" and
# "retrofit cost calculation") the left/right window placement used for the
# BC-LEEP geometry - everything downstream recalculates automatically.

$wb = $excel.ActiveWorkbook

$wsOptions = $wb.Worksheets.Item("Geometry options")
$wsCalc    = $wb.Worksheets.Item("geometry calculation")

# --- the actual data edit -------------------------------------------------
# Row 10 = "BC-LEEP" archetype, column Q = "Basement wall height" (Opt_Bsm_Height)
$wsOptions.Range("Q10").Value = 2.4384000000000001

# --- cosmetic follow-up matching how the workbook was left after editing --
# widen column C a touch on the calculation sheet
$wsCalc.Columns.Item(3).ColumnWidth = 10.5840000000000001

# move the active tab / selections to where the editor left them:
# "Geometry options" keeps a selection on the edited cell, and
# "geometry calculation" becomes the active sheet, scrolled/selected at C4.
$wsOptions.Select()
$wsOptions.Range("Q10").Select()

$wsCalc.Select()
$wsCalc.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
